# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback run:
#   - 32b3a9a4-b645-4fd3-a04e-29f6585387f3 has now been handed back (it was
#     previously "Ready for handoff"), so its status becomes "Handed back:
#     in sync with en-US" just like 8e8a92e8-51d4-4fd0-9352-55368de47472.
#   - the report lists the files in a new order: 32b3a9a4... now comes
#     first, 8e8a92e8... second.
#   - the "Latest Handback DateTime" columns on the per-locale sheets are
#     refreshed to reflect the new handback pass.
#
# This script re-writes the three worksheets (Overview, zh-cn, de-de) cell
# by cell to the post-handback state and rebuilds the hyperlinks so the
# displayed link text follows the new row order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"

$ov.Range("A3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

$ov.Range("A4").Value = ".localization-config"
$ov.Range("B4").Value = "Not to be localized"
$ov.Range("C4").Value = "Not to be localized"

# Rebuild hyperlinks in the new row order (display text follows the file
# that now occupies each row).
$ov.Cells.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/e2e/32b3a9a4-b645-4fd3-a04e-29f6585387f3.md", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/e2e/8e8a92e8-51d4-4fd0-9352-55368de47472.md", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/.localization-config", [ref]$null, [ref]$null, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-08 01:43:46"
$zh.Range("E2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md"
$zh.Range("F2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.zh-cn.xlf"
$zh.Range("G2").Value = "2016-03-08 01:44:40"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-08 01:43:46"
$zh.Range("E3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.md"
$zh.Range("F3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.zh-cn.xlf"
$zh.Range("G3").Value = "2016-03-08 01:44:40"
$zh.Range("H3").Value = "Include"

$zh.Range("A4").Value = ".localization-config"
$zh.Range("B4").Value = "Not to be localized"
$zh.Range("D4").Value = "0001-01-01 00:00:00"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Ignored"

$zh.Cells.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/e2e/32b3a9a4-b645-4fd3-a04e-29f6585387f3.md", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/794fccc07e4ad090d7b678b93a340b578cde15c1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.zh-cn.xlf", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/87109b5aa4d8f98364a1aefce2a22a748e7f72a5/e2e/32b3a9a4-b645-4fd3-a04e-29f6585387f3.md", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b8d76b3b508c17551055302f995f681fff8b618e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.zh-cn.xlf", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/e2e/8e8a92e8-51d4-4fd0-9352-55368de47472.md", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/794fccc07e4ad090d7b678b93a340b578cde15c1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/mt/8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.zh-cn.xlf", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/87109b5aa4d8f98364a1aefce2a22a748e7f72a5/e2e/8e8a92e8-51d4-4fd0-9352-55368de47472.md", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b8d76b3b508c17551055302f995f681fff8b618e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.zh-cn.xlf", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/.localization-config", [ref]$null, [ref]$null, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.de-de.xlf"
$de.Range("D2").Value = "2016-03-08 01:43:56"
$de.Range("E2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md"
$de.Range("F2").Value = "32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.de-de.xlf"
$de.Range("G2").Value = "2016-03-08 01:44:59"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.de-de.xlf"
$de.Range("D3").Value = "2016-03-08 01:43:56"
$de.Range("E3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.md"
$de.Range("F3").Value = "8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.de-de.xlf"
$de.Range("G3").Value = "2016-03-08 01:44:59"
$de.Range("H3").Value = "Include"

$de.Range("A4").Value = ".localization-config"
$de.Range("B4").Value = "Not to be localized"
$de.Range("D4").Value = "0001-01-01 00:00:00"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Ignored"

$de.Cells.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/e2e/32b3a9a4-b645-4fd3-a04e-29f6585387f3.md", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c04c45253bee76e2bbe00f008007a511c0336cf9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.de-de.xlf", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f4d40accb94f773e62dcca6a846c12ecd34849ad/e2e/32b3a9a4-b645-4fd3-a04e-29f6585387f3.md", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ca3f035a383e57b7850ae7fad4464823e6db1fcd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.de-de.xlf", [ref]$null, [ref]$null, "32b3a9a4-b645-4fd3-a04e-29f6585387f3.1a1d2eca47573bcc5b46bcf954966e79e7d05255.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/e2e/8e8a92e8-51d4-4fd0-9352-55368de47472.md", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c04c45253bee76e2bbe00f008007a511c0336cf9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/mt/8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.de-de.xlf", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/f4d40accb94f773e62dcca6a846c12ecd34849ad/e2e/8e8a92e8-51d4-4fd0-9352-55368de47472.md", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.md") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ca3f035a383e57b7850ae7fad4464823e6db1fcd/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.de-de.xlf", [ref]$null, [ref]$null, "8e8a92e8-51d4-4fd0-9352-55368de47472.c752fa8f6b6a0cdbf9ee11e410e42734798c3289.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/32c30795f540892b9ccd8ec2391c52e20b4278c2/.localization-config", [ref]$null, [ref]$null, ".localization-config") | Out-Null

Write-Output "Report regenerated for handback."
